$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 2, shifting existing data rows (2-5) down to (3-6)
$ws.Rows.Item(2).Insert()

# Populate the newly inserted row 2 with the new parameter values
$newRow = @(
    6.9531536356936074,
    3.0292758261107453,
    -9.8920303965864154,
    3.9314550243527413,
    447.46645908506372,
    14.321797140994221,
    501.58152300241562,
    9.7108031156843939,
    0.10234359962463384,
    0.10234359962463384,
    9.5731655186720435,
    1.3675950740960061,
    0.72372358143811788,
    -0.5141729627567454
)

for ($col = 1; $col -le $newRow.Length; $col++) {
    $ws.Cells.Item(2, $col).Value = $newRow[$col - 1]
}

# Update the remaining (shifted) rows' values to match the new recalculation
$rowData = @{
    3 = @(-11.700210064007177, 6.3079014360083114, -39.716940310561888, 8.5373458540983318, 581.09079153160906, 35.20362285781345, 521.82070652978575, 14.391340796572425, 0.00012210381578882254, 0.00012210381578882254, 29.082456476115009, 4.1546366394450009, -0.99966548145264511, -0.20308141008130567)
    4 = @(69.443739316870108, 5.3452346388288836, -31.057603325474204, 5.8253760186883277, 206.39405727167699, 27.870524717388513, 625.87789957676966, 8.3636593425848478, 0.12852181037151866, 0.12852181037151866, 26.56336310572247, 3.7947661579603529, 0.99365226465751344, -0.4565034527184329)
    5 = @(59.765553578621372, 10.444333552757229, -6.6506159102473026, 9.9544160369388059, 296.59722964402755, 47.08842895911711, 548.46424727478143, 16.180735901770582, 0.092516911764747425, 0.092516911764747425, 30.39501235223366, 4.342144621747666, 0.97390442398624621, -0.2987198049993256)
    6 = @(28.605119600408088, 3.3269238562436927, -5.4726940753706232, 4.4574508795781949, 411.71364427308748, 17.965757382521478, 546.02421888226456, 8.5243247649291654, 0.13960283655242786, 0.13960283655242786, 8.5825075018291592, 1.2260725002613084, 0.98379559669538819, -0.3720243894333069)
}

foreach ($r in $rowData.Keys) {
    $vals = $rowData[$r]
    for ($col = 1; $col -le $vals.Length; $col++) {
        $ws.Cells.Item($r, $col).Value = $vals[$col - 1]
    }
}
